# Daily attendance processing - normalize "Recorded By" (column G) entries so
# that any "System" / "system" token is moved to the end of the
# comma-separated list, with lowercase "system" ordered before the
# capitalized "System" when both appear together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-RecordedBy {
    param([string]$value)

    if ([string]::IsNullOrEmpty($value)) {
        return $value
    }

    $parts = $value -split ", "
    if ($parts.Count -le 1) {
        return $value
    }

    $systemParts = New-Object System.Collections.ArrayList
    $otherParts = New-Object System.Collections.ArrayList
    foreach ($p in $parts) {
        if ($p.ToLower().Equals("system")) {
            [void]$systemParts.Add($p)
        } else {
            [void]$otherParts.Add($p)
        }
    }

    if ($systemParts.Count -eq 0) {
        return $value
    }

    # Lowercase "system" sorts before capitalized "System" when tied.
    $sortedSystem = $systemParts | Sort-Object -Descending

    $result = New-Object System.Collections.ArrayList
    foreach ($p in $otherParts) { [void]$result.Add($p) }
    foreach ($p in $sortedSystem) { [void]$result.Add($p) }

    return ($result -join ", ")
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text
    if ([string]::IsNullOrEmpty($current)) {
        continue
    }
    $updated = Transform-RecordedBy $current
    if (-not $updated.Equals($current)) {
        $cell.Value = $updated
    }
}
